$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.06663987991766124
$ws.Range("C2").Value = 0.998782647210843
$ws.Range("D2").Value = 0.2049441557218043
$ws.Range("G2").Value = 0.129829331083359
$ws.Range("H2").Value = 0.99
